$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 887 so the existing weekly records (old rows
# 887-990) shift down one row (to 888-991) instead of being overwritten.
$ws.Rows("887:887").Insert()

# Populate the newly-inserted row 887 with this week's record.
$ws.Range("A887").Value = 8
$ws.Range("B887").Value = "Terminal La Palmera de La Serena"
$ws.Range("C887").Value = "Coquimbo"
$ws.Range("D887").Value = 45194
$ws.Range("E887").Value = 4
$ws.Range("F887").Value = 100112045
$ws.Range("G887").Value = "Zapallo"
$ws.Range("H887").Value = "Camote"
$ws.Range("I887").Value = "1a nueva(o)"
$ws.Range("J887").Value = 1000
$ws.Range("K887").Value = 1000
$ws.Range("L887").Value = 1100
$ws.Range("M887").Value = 1050
$ws.Range("N887").Value = "$/kilo (volumen en unidades)"
$ws.Range("O887").Value = "Perú"
$ws.Range("P887").Value = 1050
$ws.Range("Q887").Value = 1
$ws.Range("R887").Value = "Hortaliza"
